$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force a literal text value even when it looks like a date/number,
    # then drop back to the default style so no extra formatting sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------- Row 3 ----------------
$ws.Range("A3").Value = 131107108
$ws.Range("B3").Value = 79000
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("P3").Value = "SÖ Svarttjärnen, Mpd"
$ws.Range("Q3").Value = 600258
$ws.Range("R3").Value = 6952256
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Västernorrland"
$ws.Range("U3").Value = "Sundsvall"
$ws.Range("V3").Value = "Medelpad"
$ws.Range("W3").Value = "Liden"
$ws.Range("X3").Value = "2025_0061"
Set-TextValue $ws.Range("Y3") "2025-06-03"
$ws.Range("Z3").Value = "13:51"
Set-TextValue $ws.Range("AA3") "2025-06-03"
$ws.Range("AB3").Value = "13:51"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "David Isaksson"
$ws.Range("AX3").Value = "Erik Lagerin"
$ws.Range("AY3").Value = "Kustpaketet"

# ---------------- Row 4 ----------------
$ws.Range("A4").Value = 131107106
$ws.Range("B4").Value = 57881
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
Set-TextValue $ws.Range("I4") "1"
$ws.Range("P4").Value = "SÖ Svarttjärnen, Mpd"
$ws.Range("Q4").Value = 600259
$ws.Range("R4").Value = 6952247
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Västernorrland"
$ws.Range("U4").Value = "Sundsvall"
$ws.Range("V4").Value = "Medelpad"
$ws.Range("W4").Value = "Liden"
$ws.Range("X4").Value = "2025_0063"
Set-TextValue $ws.Range("Y4") "2025-06-03"
$ws.Range("Z4").Value = "13:53"
Set-TextValue $ws.Range("AA4") "2025-06-03"
$ws.Range("AB4").Value = "13:53"
$ws.Range("AC4").Value = "Gammalt födosök av spillkråka"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "David Isaksson"
$ws.Range("AX4").Value = "Samuel Koont"
$ws.Range("AY4").Value = "Kustpaketet"

# ---------------- Row 5 ----------------
$ws.Range("A5").Value = 131107103
$ws.Range("B5").Value = 8440
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 106554
$ws.Range("F5").Value = "Björksplintborre"
$ws.Range("G5").Value = "Scolytus ratzeburgii"
$ws.Range("H5").Value = "Janson, 1856"
$ws.Range("P5").Value = "SÖ Svarttjärnen, Mpd"
$ws.Range("Q5").Value = 600292
$ws.Range("R5").Value = 6952283
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Västernorrland"
$ws.Range("U5").Value = "Sundsvall"
$ws.Range("V5").Value = "Medelpad"
$ws.Range("W5").Value = "Liden"
$ws.Range("X5").Value = "2025_0066"
Set-TextValue $ws.Range("Y5") "2025-06-03"
$ws.Range("Z5").Value = "13:57"
Set-TextValue $ws.Range("AA5") "2025-06-03"
$ws.Range("AB5").Value = "13:57"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "David Isaksson"
$ws.Range("AX5").Value = "Erik Lagerin"
$ws.Range("AY5").Value = "Kustpaketet"
